$wb = $excel.ActiveWorkbook

# --- Rename existing Sheet1 to TC002 ---
$ws1 = $wb.ActiveSheet
$ws1.Name = "TC002"

# --- Add TC003 right after TC002 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TC003"

$ws2.Range("A1").Value = "MenuName"
$ws2.Range("B1").Value = "Product Id"
$ws2.Range("C1").Value = "Quantity"
$ws2.Range("D1").Value = "Size"
$ws2.Range("E1").Value = "Color"

$ws2.Range("A2").Value = "Dresses"
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = "L"
$ws2.Range("E2").Value = "Blue"

$ws2.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 18.333333333333332
$ws2.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(5).ColumnWidth = 19.666666666666668

[void]$ws2.Range("A1:E2").Select()

# --- Add TC004 right after TC003 ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "TC004"

$ws3.Range("A1").Value = "MenuName"
$ws3.Range("B1").Value = "Product Id"
$ws3.Range("C1").Value = "Quantity"
$ws3.Range("D1").Value = "Size"
$ws3.Range("E1").Value = "Color"

$ws3.Range("A2").Value = "Dresses"
$ws3.Range("B2").Value = 5
$ws3.Range("C2").Value = 3
$ws3.Range("D2").Value = "L"
$ws3.Range("E2").Value = "Blue"

$ws3.Range("A3").Value = "Women"
$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = 5
$ws3.Range("D3").Value = "M"
$ws3.Range("E3").Value = "White"
$ws3.Range("E3").Font.Name = "Consolas"
$ws3.Range("E3").Font.Size = 9
$ws3.Range("E3").Font.Color = 2236962
$ws3.Range("E3").Font.Family = 3

$ws3.Range("A4").Value = "T-shirts"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 2
$ws3.Range("D4").Value = "S"
$ws3.Range("E4").Value = "Orange"

$ws3.Range("A6").Value = "Expected Count"
$ws3.Range("B6").Value = 3

$ws3.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws3.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws3.Columns.Item(3).ColumnWidth = 13.0
$ws3.Columns.Item(4).ColumnWidth = 11.5
$ws3.Columns.Item(5).ColumnWidth = 11.333333333333334

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

[void]$ws3.Range("B6").Select()
$ws3.Activate()
